$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set Pass/Fail column (F) to "PASS" for the completed test case rows
$ws.Range("F2").Value = "PASS"
$ws.Range("F4").Value = "PASS"
$ws.Range("F5").Value = "PASS"
$ws.Range("F7").Value = "PASS"
$ws.Range("F9").Value = "PASS"
$ws.Range("F15").Value = "PASS"
$ws.Range("F17").Value = "PASS"

# E19 gets the closing message text (same as D19)
$ws.Range("E19").Value = "Thanks for shipping with Seneca Polytechnic!"

# F15/F17/F19 pick up the centered + wrapped style used elsewhere in column D
$ws.Range("F15").WrapText = $true
$ws.Range("F15").VerticalAlignment = -4108
$ws.Range("F17").WrapText = $true
$ws.Range("F17").VerticalAlignment = -4108
$ws.Range("F19").WrapText = $true
$ws.Range("F19").VerticalAlignment = -4108

# Update selection to F11
$ws.Range("F11").Select()
